$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("everything")

$ws.Rows.Item(436).Insert()
$ws.Range("A436").Value2 = "![unicorn-moonwalk](/images/unicorn-moonwalk.gif)"

$ws.Rows.Item(387).Insert()
$ws.Range("A387").Value2 = "![spyparrot](/images/spyparrot.gif)"

$ws.Rows.Item(349).Insert()
$ws.Range("A349").Value2 = "![revolutionparrot](/images/revolutionparrot.gif)"

$ws.Rows.Item(318).Insert()
$ws.Range("A318").Value2 = "![prideparrot](/images/prideparrot.gif)"

$ws.Rows.Item(278).Insert()
$ws.Range("A278").Value2 = "![partychewbacca](/images/partychewbacca.gif)"

$ws.Rows.Item(260).Insert()
$ws.Range("A260").Value2 = "![party_dumpster_fire](/images/party_dumpster_fire.gif)"

$ws.Rows.Item(146).Insert()
$ws.Range("A146").Value2 = "![gritty-look](/images/gritty-look.gif)"

$ws.Rows.Item(87).Insert()
$ws.Range("A87").Value2 = "![dancing-unicorn](/images/dancing-unicorn.gif)"

$ws.Application.GoTo($ws.Range("D463"), $true)

Write-Host "done"
